$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Implementacion de intefaz punto
$ws.Range("A6").Value = "Implementacion de intefaz punto"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 0.0048611111111111112
$ws.Range("E6").Value = 0.34166666666666662
$ws.Range("F6").Value = 0.34722222222222227
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# Row 7 - Implementacion de Punto2D
$ws.Range("A7").Value = "Implementacion de Punto2D"
$ws.Range("B7").Value = 80
$ws.Range("C7").Value = 107
$ws.Range("D7").Value = 0.013888888888888888
$ws.Range("E7").Value = 0.34791666666666665
$ws.Range("F7").Value = 0.36041666666666666
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# Row 8 - Implementacion de Punto3D
$ws.Range("A8").Value = "Implementacion de Punto3D"
$ws.Range("B8").Value = 80
$ws.Range("D8").Value = 0.013888888888888888
$ws.Range("E8").Value = 0.36458333333333331
$ws.Range("F8").Value = 0.38194444444444442
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

$excel.CalculateFullRebuild()

# A few formula cells (previously cached as text "-") don't get their
# result type/value refreshed by a plain recalculation, so re-apply their
# formulas to force Excel to re-evaluate and re-type them.
$ws.Range("B12").Formula = $ws.Range("B12").Formula()
$ws.Range("B16").Formula = $ws.Range("B16").Formula()
$ws.Range("B17").Formula = $ws.Range("B17").Formula()
$ws.Range("B18").Formula = $ws.Range("B18").Formula()

$excel.CalculateFullRebuild()

$ws.Range("A9").Select()

$wb.Save()
